# Update the ELO leaderboard: re-rank players by their new ELO score and
# drop the players who no longer appear (Jean-Pierre, Gergo, Kirill,
# Patricia, Nic, Peter, Gabriel). Row 1 (header) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New standings, already sorted by elo (desc), row 2 downward.
$data = @(
    ,("Eric", 1874.7)
    ,("Jerome", 1865.5)
    ,("Sean", 1760.2)
    ,("Saravanan", 1723.1)
    ,("Hassan", 1721.3)
    ,("Paulo", 1701.4)
    ,("Alexis", 1691.1)
    ,("Anil", 1677.1)
    ,("Zhengnan", 1675.4)
    ,("Ayman", 1653.9)
    ,("Xi", 1637.8)
    ,("Sadeed", 1631.4)
    ,("Mark", 1630.9)
    ,("Eugene", 1611.8)
    ,("Charlie", 1605.9)
    ,("Luis", 1555.8)
    ,("Karla", 1525.8)
    ,("Fernando", 1521.7)
    ,("Dylan", 1521.2)
    ,("Jameel", 1516)
    ,("Jofrey", 1515.8)
    ,("Yevhen", 1513.7)
    ,("Madi", 1499.4)
    ,("Octavio", 1465.2)
    ,("Abdurauf", 1449.6)
    ,("Amanat", 1430.1)
    ,("Omar", 1410.6)
    ,("Carlos", 1405.2)
    ,("Rawan", 1400.8)
    ,("Faruk", 1381.6)
    ,("Abdulmajeed", 1375.7)
    ,("Mustafa", 1372.3)
    ,("Lucas", 1325)
    ,("Haytham", 1247.5)
    ,("Huawen", 1231.4)
    ,("Xingzhu", 1208.4)
    ,("Hashim", 1206.3)
    ,("Juris", 1198.7)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}

# The old list had 45 data rows (to row 46); the new one only has 38
# (to row 39). Wipe the now-unused tail so the used range shrinks back
# to A1:B39.
$ws.Range("A40:B46").ClearContents()
